$wb = $excel.ActiveWorkbook

# Sheets that contain the "想去人数" (interested count) data that changed: 展览 and 全部类型
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 82
    $ws.Range("F4").Value = 1476
    $ws.Range("F6").Value = 31
    $ws.Range("F7").Value = 116
    $ws.Range("F9").Value = 263
}

$wb.Save()
